$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Range("H94").Value = 2454.2
$ws.Range("I94").Value = 2454.2
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2454.2
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2003.2
$ws.Range("N94").ClearContents()
# Row 132
$ws.Range("H132").Value = 938.8421
$ws.Range("I132").Value = 933.6286
$ws.Range("J132").Value = 999.6667
$ws.Range("K132").Value = 2800.8858
$ws.Range("L132").Value = 2999.0001
$ws.Range("M132").Value = -270.8858
$ws.Range("N132").Value = -8059.0001
# Row 138
$ws.Range("H138").Value = 1762.4138
$ws.Range("J138").Value = 2224.875
$ws.Range("L138").Value = 6674.625
$ws.Range("N138").Value = -16954.625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4252.383
$ws.Range("I32").Value = 2281.3057
$ws.Range("J32").Value = 10703.182
$ws.Range("K32").Value = 2281.3057
$ws.Range("L32").Value = 10703.182
$ws.Range("M32").Value = -1994.3057
$ws.Range("N32").Value = -11277.182
# Row 61
$ws.Range("H61").Value = 3400.85
$ws.Range("I61").Value = 2392.4707
$ws.Range("K61").Value = 2392.4707
$ws.Range("M61").Value = -2180.4707
# Row 74
$ws.Range("H74").Value = 725.6
$ws.Range("I74").Value = 725.6
$ws.Range("K74").Value = 725.6
$ws.Range("M74").Value = 148.4
# Row 77
$ws.Range("H77").Value = 725.6
$ws.Range("I77").Value = 725.6
$ws.Range("K77").Value = 3628
$ws.Range("M77").Value = 740
# Row 109
$ws.Range("H109").Value = 42496
$ws.Range("J109").Value = 42496
$ws.Range("L109").Value = 42496
$ws.Range("N109").Value = -45270
# Row 132
$ws.Range("H132").Value = 1348.9032
$ws.Range("I132").Value = 1068.5
$ws.Range("K132").Value = 3205.5
$ws.Range("M132").Value = -675.5
# Row 136
$ws.Range("H136").Value = 3400.85
$ws.Range("I136").Value = 2392.4707
$ws.Range("K136").Value = 7177.4121
$ws.Range("M136").Value = -4627.4121

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2288.9167
$ws.Range("I20").Value = 2196
$ws.Range("J20").Value = 3311
$ws.Range("K20").Value = 2196
$ws.Range("L20").Value = 3311
$ws.Range("M20").Value = -1949
$ws.Range("N20").Value = -3805
# Row 99
$ws.Range("H99").Value = 1059.6
$ws.Range("I99").Value = 999.5
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 999.5
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = 498.5
$ws.Range("N99").Value = -4296
# Row 105
$ws.Range("H105").Value = 2093.75
$ws.Range("I105").Value = 2228.5715
$ws.Range("K105").Value = 2228.5715
$ws.Range("M105").Value = -481.5715
# Row 134
$ws.Range("H134").Value = 10834.061
$ws.Range("I134").Value = 11041.148
$ws.Range("J134").Value = 9902.166999999999
$ws.Range("K134").Value = 33123.444
$ws.Range("L134").Value = 29706.501
$ws.Range("M134").Value = -30588.444
$ws.Range("N134").Value = -34776.501

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2777.1667
$ws.Range("I31").Value = 1524.9412
$ws.Range("J31").Value = 5818.2856
$ws.Range("K31").Value = 1524.9412
$ws.Range("L31").Value = 5818.2856
$ws.Range("M31").Value = -1229.9412
$ws.Range("N31").Value = -6408.2856
# Row 34
$ws.Range("H34").Value = 2777.1667
$ws.Range("I34").Value = 1524.9412
$ws.Range("J34").Value = 5818.2856
$ws.Range("K34").Value = 1524.9412
$ws.Range("L34").Value = 5818.2856
$ws.Range("M34").Value = -1322.9412
$ws.Range("N34").Value = -6222.2856
# Row 58
$ws.Range("H58").Value = 1977740.9
$ws.Range("I58").Value = 3953726.5
$ws.Range("K58").Value = 3953726.5
$ws.Range("M58").Value = -3953523.5
# Row 74
$ws.Range("H74").Value = 27375
$ws.Range("J74").Value = 27375
$ws.Range("L74").Value = 27375
$ws.Range("N74").Value = -29123
# Row 77
$ws.Range("H77").Value = 27375
$ws.Range("J77").Value = 27375
$ws.Range("L77").Value = 82125
$ws.Range("N77").Value = -90861
# Row 132
$ws.Range("H132").Value = 1366.6285
$ws.Range("I132").Value = 869.5
$ws.Range("J132").Value = 2802.7778
$ws.Range("K132").Value = 2608.5
$ws.Range("L132").Value = 8408.3334
$ws.Range("M132").Value = -78.5
$ws.Range("N132").Value = -13468.3334
# Row 136
$ws.Range("H136").Value = 1977740.9
$ws.Range("I136").Value = 3953726.5
$ws.Range("K136").Value = 11861179.5
$ws.Range("M136").Value = -11858629.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 61
$ws.Range("H61").Value = 200
$ws.Range("I61").Value = 200
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 600
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -385
$ws.Range("N61").ClearContents()
# Row 107
$ws.Range("H107").Value = 697.75
$ws.Range("I107").Value = 290
$ws.Range("J107").Value = 719.2105
$ws.Range("K107").Value = 870
$ws.Range("L107").Value = 2157.6315
$ws.Range("M107").Value = 1050
$ws.Range("N107").Value = -5997.6315
# Row 113
$ws.Range("H113").Value = 8607.385
$ws.Range("J113").Value = 1029.2
$ws.Range("L113").Value = 3087.6
$ws.Range("N113").Value = -7427.6
# Row 122
$ws.Range("H122").Value = 764.8461
$ws.Range("I122").Value = 591.2
$ws.Range("J122").Value = 873.375
$ws.Range("K122").Value = 5320.8
$ws.Range("L122").Value = 7860.375
$ws.Range("M122").Value = -2870.8
$ws.Range("N122").Value = -12760.375
# Row 131
$ws.Range("H131").Value = 11751.094
$ws.Range("J131").Value = 12499.033
$ws.Range("L131").Value = 37497.099
$ws.Range("N131").Value = -47577.099

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2489.1738
$ws.Range("I102").Value = 2579.75
$ws.Range("J102").Value = 2390.3635
$ws.Range("K102").Value = 2579.75
$ws.Range("L102").Value = 2390.3635
$ws.Range("M102").Value = -957.75
$ws.Range("N102").Value = -5634.363499999999
# Row 132
$ws.Range("H132").Value = 1427844.5
$ws.Range("I132").Value = 2026108.4
$ws.Range("J132").Value = 6967.625
$ws.Range("K132").Value = 6078325.199999999
$ws.Range("L132").Value = 20902.875
$ws.Range("M132").Value = -6075795.199999999
$ws.Range("N132").Value = -25962.875

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3349.9167
$ws.Range("I16").Value = 3714.8
$ws.Range("J16").Value = 1525.5
$ws.Range("K16").Value = 3714.8
$ws.Range("L16").Value = 1525.5
$ws.Range("M16").Value = -3544.8
$ws.Range("N16").Value = -1865.5
# Row 46
$ws.Range("H46").Value = 1451.6666
$ws.Range("I46").Value = 210
$ws.Range("K46").Value = 210
$ws.Range("M46").Value = -22
# Row 61
$ws.Range("H61").Value = 2771.6667
$ws.Range("I61").Value = 2810.8
$ws.Range("J61").Value = 2693.4
$ws.Range("K61").Value = 2810.8
$ws.Range("L61").Value = 2693.4
$ws.Range("M61").Value = -2608.8
$ws.Range("N61").Value = -3097.4
# Row 113
$ws.Range("H113").Value = 2771.6667
$ws.Range("I113").Value = 2810.8
$ws.Range("J113").Value = 2693.4
$ws.Range("K113").Value = 2810.8
$ws.Range("L113").Value = 2693.4
$ws.Range("M113").Value = -640.8000000000002
$ws.Range("N113").Value = -7033.4
# Row 136
$ws.Range("H136").Value = 4089.2104
$ws.Range("J136").Value = 5869.4
$ws.Range("L136").Value = 17608.2
$ws.Range("N136").Value = -22708.2

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 2043.7142
$ws.Range("I14").Value = 700
$ws.Range("J14").Value = 2267.6667
$ws.Range("K14").Value = 700
$ws.Range("L14").Value = 2267.6667
$ws.Range("M14").Value = -532
$ws.Range("N14").Value = -2603.6667
# Row 100
$ws.Range("H100").Value = 502
$ws.Range("I100").Value = 304
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 608
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -67
$ws.Range("N100").Value = -2482
# Row 126
$ws.Range("H126").Value = 5001.385
$ws.Range("I126").Value = 4379.9644
$ws.Range("J126").Value = 6583.1816
$ws.Range("K126").Value = 13139.8932
$ws.Range("L126").Value = 19749.5448
$ws.Range("M126").Value = -10669.8932
$ws.Range("N126").Value = -24689.5448
# Row 132
$ws.Range("H132").Value = 6417.3335
$ws.Range("I132").Value = 1252.2727
$ws.Range("K132").Value = 3756.8181
$ws.Range("M132").Value = -1226.8181
# Row 136
$ws.Range("H136").Value = 27780940
$ws.Range("I136").Value = 39685744
$ws.Range("K136").Value = 119057232
$ws.Range("M136").Value = -119054682
